$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 14: Status "in Arbeit" -> "fertig", add Akzeptanztest "akzeptiert" ---
$ws.Cells.Item(14, 2).Value = "fertig"
$ws.Cells.Item(14, 10).Value = "akzeptiert"

# --- Row 20: Status "jungfräulich" -> "fertig", add Akzeptanztest "akzeptiert" ---
$ws.Cells.Item(20, 2).Value = "fertig"
$ws.Cells.Item(20, 10).Value = "akzeptiert"

# --- Row 21: add Akzeptanztest "akzeptiert" ---
$ws.Cells.Item(21, 10).Value = "akzeptiert"

# --- Row 22: add Akzeptanztest "akzeptiert" ---
$ws.Cells.Item(22, 10).Value = "akzeptiert"

# --- Row 23: Status "in Arbeit" -> "fertig", add Akzeptanztest, actual effort & completion date ---
$ws.Cells.Item(23, 2).Value = "fertig"
$ws.Cells.Item(23, 10).Value = "akzeptiert"
$ws.Cells.Item(23, 12).Value = "2h 30min"
# copy the date format (s="3") from an existing date cell (E23) before writing the value
$ws.Cells.Item(23, 5).Copy()
$ws.Cells.Item(23, 13).PasteSpecial(-4122)
$ws.Cells.Item(23, 13).Value = 40830

# --- New rows 24 ("Vortrag") & 25 ("Dokumentation") StoryCards ---
# Shared-string table is append-on-first-use, so the cells below are written in
# the specific order needed to reproduce the target string indices 80-85:
#   80 "2h 30min" (written above, L23), 81 "4h", 82 "Vortrag", 83 "Dokumentation",
#   84 "Präsentation erstellen", 85 "Dokumentation erstellen"
$ws.Cells.Item(24, 1).Value = 17
$ws.Cells.Item(24, 2).Value = "in Arbeit"
$ws.Cells.Item(24, 3).Value = "hoch"
$ws.Cells.Item(24, 6).Value = "Wiederschein"
$ws.Cells.Item(24, 7).Value = "alle"
$ws.Cells.Item(24, 8).Value = "Ausarbeitung"
$ws.Cells.Item(24, 11).Value = "4h"
$ws.Cells.Item(24, 4).Value = "Vortrag"

$ws.Cells.Item(25, 1).Value = 18
$ws.Cells.Item(25, 2).Value = "in Arbeit"
$ws.Cells.Item(25, 3).Value = "hoch"
$ws.Cells.Item(25, 6).Value = "Wiederschein"
$ws.Cells.Item(25, 7).Value = "alle"
$ws.Cells.Item(25, 8).Value = "Ausarbeitung"
$ws.Cells.Item(25, 4).Value = "Dokumentation"

$ws.Cells.Item(24, 9).Value = "Präsentation erstellen"
$ws.Cells.Item(25, 9).Value = "Dokumentation erstellen"
$ws.Cells.Item(25, 11).Value = "8h"

# copy the date format (s="3") from an existing date cell before writing values
$ws.Cells.Item(23, 5).Copy()
$ws.Cells.Item(24, 5).PasteSpecial(-4122)
$ws.Cells.Item(24, 5).Value = 40840

$ws.Cells.Item(23, 5).Copy()
$ws.Cells.Item(25, 5).PasteSpecial(-4122)
$ws.Cells.Item(25, 5).Value = 40840

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 33.75
$ws.Rows.Item(7).RowHeight = 25.5

# --- Column I width (no longer best-fit, slightly wider) ---
$ws.Columns.Item(9).ColumnWidth = 38.57

# --- Selection / view: move to H27, drop the frozen/topLeft A4 scroll position ---
$ws.Range("H27").Select()
